$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6955972909927368
$ws.Range("B1").Value = 0.9964777231216431
$ws.Range("C1").Value = 1.901272177696228
$ws.Range("D1").Value = 3.309988975524902
$ws.Range("E1").Value = 3.64540433883667
